# Drop in RMI script files
# - Remove the "Data Texas" worksheet (local analyst scratch sheet / notes)
# - Restore the canonical "HPPECbP" excess-capacity value of 25% (0.25)

$wb = $excel.ActiveWorkbook

# Delete the "Data Texas" worksheet entirely (suppress the standard
# "permanently delete this sheet" confirmation prompt).
$excel.DisplayAlerts = $false
$dataTexas = $wb.Worksheets.Item("Data Texas")
$dataTexas.Delete()
$excel.DisplayAlerts = $true

# Restore the original excess-capacity assumption on the HPPECbP sheet.
$ws = $wb.Worksheets.Item("HPPECbP")
$ws.Range("B2").Value = 0.25

# Make "About" the active sheet again (matches the pristine, freshly
# dropped-in copy of the source file).
$wb.Worksheets.Item("About").Activate()
